$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen description / acceptance-criteria columns (E, F) ---
$ws.Columns.Item(5).ColumnWidth = 69.0
$ws.Columns.Item(6).ColumnWidth = 75.5

# --- Freeze the header row and scroll the view back to the top ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- Append new backlog item LALA-027 as row 28 ---
$ws.Range("A28").Value = "LALA-027"
$ws.Range("B28").Value = "P1"
$ws.Range("C28").Value = "UI  "
$ws.Range("D28").Value = "Tire change time changing total fuel unecessarily"
$ws.Range("E28").Value = "Changing tire time that is still less than refuel time is adjusting fuel when it shouldn't because the time is already accounted by the refuel time."
$ws.Range("F28").Value = "When tire change time is less than refuel time, the total fuel needed should not change"
$ws.Range("G28").Value = "Backlog"
$ws.Range("H28").Value = "Andy"
$ws.Range("I28").Value = "any"
$ws.Rows.Item(28).RowHeight = 30

# --- Final selection lands on the new row's last populated cell ---
$ws.Range("J28").Select() | Out-Null
